# Add start/stop for Akeelah
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("film_info")

# Row 2 is the "Akeelah" entry (A2 = Akeelah, B2 = Akeelah and the Bee (2006)).
# Columns C/D (Clip_Start / Clip_Stop) are formatted as Text (style index 3,
# inherited from the column style) so assign the time strings as text values.
$ws.Cells.Item(2, 3).NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "01:42:23"

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "01:45:07"

$ws.Range("C3").Select()
